# Update crypto price/volume figures to the latest scrape (GitHub Actions run).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row -> column letter -> new cell text.
$updates = @{
    2 = @{ 'D' = '60.239.78'; 'E' = '  +1.92%  ' }
    3 = @{ 'D' = '2.680.70' }
    4 = @{ 'E' = '  +0.06%  ' }
    5 = @{ 'D' = '522.89'; 'E' = '  +0.90%  ' }
    6 = @{ 'D' = '146.55'; 'E' = '  +0.48%  ' }
    7 = @{ 'D' = '0.998'; 'E' = '  +0.34%  ' }
    8 = @{ 'D' = '0.576'; 'E' = '  +1.28%  ' }
    9 = @{ 'D' = '2.698.66'; 'E' = '  -0.77%  ' }
    10 = @{ 'D' = '6.48'; 'E' = '  +3.36%  ' }
    11 = @{ 'D' = '0.106'; 'E' = '  -1.14%  ' }
    12 = @{ 'E' = '  +0.50%  ' }
    13 = @{ 'E' = '  +1.62%  ' }
    14 = @{ 'D' = '3.160.95'; 'E' = '  +0.06%  ' }
    15 = @{ 'D' = '60.429.27'; 'E' = '  +2.30%  ' }
    16 = @{ 'D' = '21.36'; 'E' = '  +0.87%  ' }
    17 = @{ 'E' = '  -0.07%  ' }
    18 = @{ 'D' = '2.694.87'; 'E' = '  -0.89%  ' }
    19 = @{ 'D' = '351.35'; 'E' = '  +1.30%  ' }
    20 = @{ 'D' = '4.54'; 'E' = '  -0.63%  ' }
    21 = @{ 'D' = '10.55'; 'E' = '  +0.38%  ' }
    22 = @{ 'D' = '6.34'; 'E' = '  +1.94%  ' }
    23 = @{ 'E' = '  +0.07%  ' }
    24 = @{ 'D' = '63.22'; 'E' = '  +3.41%  ' }
    25 = @{ 'D' = '0.422' }
    26 = @{ 'E' = '  +4.56%  ' }
    27 = @{ 'D' = '0.994'; 'E' = '  +0.18%  ' }
    28 = @{ 'B' = 'InternetComputer(DFINITY)'; 'C' = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'; 'D' = '7.36'; 'E' = '  +1.23%  ' }
    29 = @{ 'B' = 'PEPE'; 'C' = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'; 'D' = '0.0₃0819'; 'E' = '  -0.52%  ' }
    30 = @{ 'D' = '6.82'; 'E' = '  +5.93%  ' }
    31 = @{ 'E' = '  +0.17%  ' }
    32 = @{ 'E' = '  +0.74%  ' }
    33 = @{ 'D' = '19.13'; 'E' = '  -0.08%  ' }
    34 = @{ 'D' = '147.27'; 'E' = '  -1.94%  ' }
    35 = @{ 'E' = '  +5.74%  ' }
    36 = @{ 'E' = '  +8.19%  ' }
    37 = @{ 'D' = '0.954'; 'E' = '  -6.99%  ' }
    38 = @{ 'D' = '0.877'; 'E' = '  +2.26%  ' }
    39 = @{ 'E' = '  +7.53%  ' }
    40 = @{ 'E' = '  -0.01%  ' }
    41 = @{ 'D' = '3.71'; 'E' = '  -0.43%  ' }
    42 = @{ 'D' = '284.87'; 'E' = '  +0.90%  ' }
    43 = @{ 'D' = '20.10'; 'E' = '  -1.03%  ' }
    44 = @{ 'D' = '0.0993'; 'E' = '  +0.83%  ' }
    45 = @{ 'B' = 'FirstDigitalUSD'; 'C' = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'; 'D' = '0.997'; 'E' = '  +0.45%  ' }
    46 = @{ 'B' = 'Mantle'; 'C' = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'; 'D' = '0.611'; 'E' = '  -2.50%  ' }
    47 = @{ 'D' = '2.131.98'; 'E' = '  +5.58%  ' }
    48 = @{ 'D' = '4.93'; 'E' = '  +2.98%  ' }
    49 = @{ 'E' = '  +0.98%  ' }
    50 = @{ 'E' = '  +1.58%  ' }
    51 = @{ 'D' = '19.41'; 'E' = '  +4.57%  ' }
}

# Values that look like plain numbers (e.g. "6.48") need an explicit text
# NumberFormat first so Excel keeps storing them as text instead of coercing
# them to a Double, matching how the Price column is authored upstream.
$numericPattern = '^\d+(\.\d+)?$'

foreach ($row in $updates.Keys) {
    $rowUpdates = $updates[$row]
    foreach ($col in $rowUpdates.Keys) {
        $newValue = $rowUpdates[$col]
        $range = $ws.Range("$col$row")
        if ($col -eq 'D' -and $newValue -match $numericPattern) {
            $range.NumberFormat = '@'
        }
        $range.Value = $newValue
    }
}
